$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6543
$ws.Range("C27").Value = 1016
$ws.Range("D27").Value = 6099984
$ws.Range("E27").Value = 932.2916093535075
$ws.Range("F27").Value = 9.966386554621852
$ws.Range("G27").Value = 7.286166842661035
$ws.Range("H27").Value = 25.1148649771193
